$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.132.35'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.545.23'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.41'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.06'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.546.65'
$ws.Range("E9").Value = '  +0.64%  '
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.05'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.60'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.000.54'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.096.70'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.35'
$ws.Range("E18").Value = '  +135.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.515.96'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("E20").Value = '  +4.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.07'
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '372.73'
$ws.Range("E22").Value = '  +4.61%  '
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.82'
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.93'
$ws.Range("E27").Value = '  -4.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.654.36'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0975'
$ws.Range("E30").Value = '  -0.93%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '543.37'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.44'
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.69'
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.22'
$ws.Range("E39").Value = '  +2.85%  '
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("E44").Value = '  +2.30%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.35'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0288'
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '148.29'
$ws.Range("E48").Value = '  -0.71%  '
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.553'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("E51").Value = '  +1.26%  '
